$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2026-02-07 -> 2026-02-08, i.e. 46060 -> 46061) for every data row (2..500).
$ws.Range("C2:C500").Value = 46061
